# Apply the data corrections to the CONTRACTS database sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (order_id 2): stone_type changed and price corrected
$ws.Range("C3").Value = "Sandstein"
$ws.Range("G3").Value = 15000

# Row 6 (order_id 5): amount and price corrected
$ws.Range("D6").Value = 100
$ws.Range("G6").Value = 7500
